# Daily attendance processing - 2025-10-18 17:40:13
# Rotate the "Recorded By" (column G) comma-separated list left by one
# position: the first entry moves to the end of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ",\s*"
        if ($parts.Count -gt 1) {
            $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
            $cell.Value2 = $rotated
        }
    }
}
